# Applies the "Cleaned up and refreshed matches. Verified scores of playoff games. Included final stats." edit
# to combined_stats.xlsx (Team Stats / Individual Stats / Team Rosters).
$wb = $excel.ActiveWorkbook

# --- Sheet1: Team Stats ---
$ws1 = $wb.Worksheets.Item("Team Stats")
$ws1.Range("A2").Value = "Elky"
$ws1.Range("B2").Value = 6
$ws1.Range("C2").Value = 0
$ws1.Range("D2").Value = 0
$ws1.Range("E2").Value = 1
$ws1.Range("F2").Value = 75.33333333333333
$ws1.Range("G2").Value = 3.5
$ws1.Range("H2").Value = 18.33333333333333
$ws1.Range("I2").Value = 28.66666666666667
$ws1.Range("J2").Value = 17.16666666666667
$ws1.Range("K2").Value = 29.5
$ws1.Range("L2").Value = 452
$ws1.Range("M2").Value = 21
$ws1.Range("N2").Value = 110
$ws1.Range("O2").Value = 172
$ws1.Range("P2").Value = 103
$ws1.Range("Q2").Value = 177

$ws1.Range("A3").Value = "The Chemists"
$ws1.Range("B3").Value = 5
$ws1.Range("C3").Value = 1
$ws1.Range("D3").Value = 0
$ws1.Range("E3").Value = 0.8333333333333334
$ws1.Range("F3").Value = 51.5
$ws1.Range("G3").Value = 2.833333333333333
$ws1.Range("H3").Value = 10
$ws1.Range("I3").Value = 20.66666666666667
$ws1.Range("J3").Value = 11.83333333333333
$ws1.Range("K3").Value = 19
$ws1.Range("L3").Value = 309
$ws1.Range("M3").Value = 17
$ws1.Range("N3").Value = 60
$ws1.Range("O3").Value = 124
$ws1.Range("P3").Value = 71
$ws1.Range("Q3").Value = 114

$ws1.Range("A4").Value = "Southern Buckeye Regional Champions"
$ws1.Range("B4").Value = 4
$ws1.Range("C4").Value = 2
$ws1.Range("D4").Value = 0
$ws1.Range("E4").Value = 0.6666666666666666
$ws1.Range("F4").Value = 50.33333333333334
$ws1.Range("G4").Value = 1.833333333333333
$ws1.Range("H4").Value = 11.66666666666667
$ws1.Range("I4").Value = 16.5
$ws1.Range("J4").Value = 10.66666666666667
$ws1.Range("K4").Value = 23.16666666666667
$ws1.Range("L4").Value = 302
$ws1.Range("M4").Value = 11
$ws1.Range("N4").Value = 70
$ws1.Range("O4").Value = 99
$ws1.Range("P4").Value = 64
$ws1.Range("Q4").Value = 139

$ws1.Range("A5").Value = "AP World Government"
$ws1.Range("B5").Value = 2
$ws1.Range("C5").Value = 4
$ws1.Range("D5").Value = 0
$ws1.Range("E5").Value = 0.3333333333333333
$ws1.Range("F5").Value = 51.16666666666666
$ws1.Range("G5").Value = 3.333333333333333
$ws1.Range("H5").Value = 10.16666666666667
$ws1.Range("I5").Value = 15.66666666666667
$ws1.Range("J5").Value = 13
$ws1.Range("K5").Value = 22.5
$ws1.Range("L5").Value = 307
$ws1.Range("M5").Value = 20
$ws1.Range("N5").Value = 61
$ws1.Range("O5").Value = 94
$ws1.Range("P5").Value = 78
$ws1.Range("Q5").Value = 135

$ws1.Range("A6").Value = "Science Bros"
$ws1.Range("B6").Value = 4
$ws1.Range("C6").Value = 2
$ws1.Range("D6").Value = 0
$ws1.Range("E6").Value = 0.6666666666666666
$ws1.Range("F6").Value = 51.33333333333334
$ws1.Range("G6").Value = 2
$ws1.Range("H6").Value = 13.66666666666667
$ws1.Range("I6").Value = 15
$ws1.Range("J6").Value = 12.83333333333333
$ws1.Range("K6").Value = 23.5
$ws1.Range("L6").Value = 308
$ws1.Range("M6").Value = 12
$ws1.Range("N6").Value = 82
$ws1.Range("O6").Value = 90
$ws1.Range("P6").Value = 77
$ws1.Range("Q6").Value = 141

$ws1.Range("A7").Value = "The Woke Left Strikes Again"
$ws1.Range("B7").Value = 2
$ws1.Range("C7").Value = 4
$ws1.Range("D7").Value = 0
$ws1.Range("E7").Value = 0.3333333333333333
$ws1.Range("F7").Value = 45.5
$ws1.Range("G7").Value = 2
$ws1.Range("H7").Value = 12.16666666666667
$ws1.Range("I7").Value = 15
$ws1.Range("J7").Value = 9.166666666666666
$ws1.Range("K7").Value = 21.33333333333333
$ws1.Range("L7").Value = 273
$ws1.Range("M7").Value = 12
$ws1.Range("N7").Value = 73
$ws1.Range("O7").Value = 90
$ws1.Range("P7").Value = 55
$ws1.Range("Q7").Value = 128

$ws1.Range("A8").Value = "H-Squared"
$ws1.Range("B8").Value = 0
$ws1.Range("C8").Value = 5
$ws1.Range("D8").Value = 0
$ws1.Range("E8").Value = 0
$ws1.Range("F8").Value = 36.8
$ws1.Range("G8").Value = 1.2
$ws1.Range("H8").Value = 10.2
$ws1.Range("I8").Value = 9.6
$ws1.Range("J8").Value = 9.8
$ws1.Range("K8").Value = 17.4
$ws1.Range("L8").Value = 184
$ws1.Range("M8").Value = 6
$ws1.Range("N8").Value = 51
$ws1.Range("O8").Value = 48
$ws1.Range("P8").Value = 49
$ws1.Range("Q8").Value = 87

$ws1.Range("A9").Value = "The 'Daley' Double"
$ws1.Range("B9").Value = 0
$ws1.Range("C9").Value = 5
$ws1.Range("D9").Value = 0
$ws1.Range("E9").Value = 0
$ws1.Range("F9").Value = 21.8
$ws1.Range("G9").Value = 0.6
$ws1.Range("H9").Value = 4.6
$ws1.Range("I9").Value = 8.8
$ws1.Range("J9").Value = 5.6
$ws1.Range("K9").Value = 7.4
$ws1.Range("L9").Value = 109
$ws1.Range("M9").Value = 3
$ws1.Range("N9").Value = 23
$ws1.Range("O9").Value = 44
$ws1.Range("P9").Value = 28
$ws1.Range("Q9").Value = 37

$ws1.Range("C16").Select()

# --- Sheet2: Individual Stats ---
$ws2 = $wb.Worksheets.Item("Individual Stats")
# Remove the trailing placeholder row (previously the "#NAME?" player row)
$ws2.Range("A17:H17").EntireRow.Delete()

$ws2.Range("A2").Value = "Will Reuter"
$ws2.Range("B2").Value = 28.83333333333333
$ws2.Range("C2").Value = 1.833333333333333
$ws2.Range("D2").Value = 11.66666666666667
$ws2.Range("E2").Value = 11
$ws2.Range("F2").Value = 70
$ws2.Range("G2").Value = 173
$ws2.Range("H2").Value = 6

$ws2.Range("A3").Value = "Abigail Friedstrom"
$ws2.Range("B3").Value = 28.66666666666667
$ws2.Range("C3").Value = 2.333333333333333
$ws2.Range("D3").Value = 10.83333333333333
$ws2.Range("E3").Value = 14
$ws2.Range("F3").Value = 65
$ws2.Range("G3").Value = 172
$ws2.Range("H3").Value = 6

$ws2.Range("A4").Value = "Scott Youngquist"
$ws2.Range("B4").Value = 23.5
$ws2.Range("C4").Value = 2.5
$ws2.Range("D4").Value = 8
$ws2.Range("E4").Value = 15
$ws2.Range("F4").Value = 48
$ws2.Range("G4").Value = 141
$ws2.Range("H4").Value = 6

$ws2.Range("A5").Value = "La Maestra (Parker Johnson)"
$ws2.Range("B5").Value = 21.83333333333333
$ws2.Range("C5").Value = 1.166666666666667
$ws2.Range("D5").Value = 9.166666666666666
$ws2.Range("E5").Value = 7
$ws2.Range("F5").Value = 55
$ws2.Range("G5").Value = 131
$ws2.Range("H5").Value = 6

$ws2.Range("A6").Value = "Ky Reckamp"
$ws2.Range("B6").Value = 18.5
$ws2.Range("C6").Value = 1.166666666666667
$ws2.Range("D6").Value = 7.5
$ws2.Range("E6").Value = 7
$ws2.Range("F6").Value = 45
$ws2.Range("G6").Value = 111
$ws2.Range("H6").Value = 6

$ws2.Range("A7").Value = "Patty (Patrick Wells)"
$ws2.Range("B7").Value = 18.5
$ws2.Range("C7").Value = 1.5
$ws2.Range("D7").Value = 7
$ws2.Range("E7").Value = 9
$ws2.Range("F7").Value = 42
$ws2.Range("G7").Value = 111
$ws2.Range("H7").Value = 6

$ws2.Range("A8").Value = "Noah Mcredmond"
$ws2.Range("B8").Value = 16.66666666666667
$ws2.Range("C8").Value = 1
$ws2.Range("D8").Value = 6.833333333333333
$ws2.Range("E8").Value = 6
$ws2.Range("F8").Value = 41
$ws2.Range("G8").Value = 100
$ws2.Range("H8").Value = 6

$ws2.Range("A9").Value = "Clete Reinberger"
$ws2.Range("B9").Value = 13.66666666666667
$ws2.Range("C9").Value = 1
$ws2.Range("D9").Value = 5.333333333333333
$ws2.Range("E9").Value = 6
$ws2.Range("F9").Value = 32
$ws2.Range("G9").Value = 82
$ws2.Range("H9").Value = 6

$ws2.Range("A10").Value = "Aidan Hook"
$ws2.Range("B10").Value = 13.6
$ws2.Range("C10").Value = 0.8
$ws2.Range("D10").Value = 5.6
$ws2.Range("E10").Value = 4
$ws2.Range("F10").Value = 28
$ws2.Range("G10").Value = 68
$ws2.Range("H10").Value = 5

$ws2.Range("A11").Value = "GHAS (Andrew Ghastine)"
$ws2.Range("B11").Value = 11.83333333333333
$ws2.Range("C11").Value = 1.833333333333333
$ws2.Range("D11").Value = 3.166666666666667
$ws2.Range("E11").Value = 11
$ws2.Range("F11").Value = 19
$ws2.Range("G11").Value = 71
$ws2.Range("H11").Value = 6

$ws2.Range("A12").Value = "El Rey de Rompecabezas (Mr. Cool)"
$ws2.Range("B12").Value = 11.5
$ws2.Range("C12").Value = 0.8333333333333334
$ws2.Range("D12").Value = 4.5
$ws2.Range("E12").Value = 5
$ws2.Range("F12").Value = 27
$ws2.Range("G12").Value = 69
$ws2.Range("H12").Value = 6

$ws2.Range("A13").Value = "Axe (Hildy)"
$ws2.Range("B13").Value = 10.4
$ws2.Range("C13").Value = 0.4
$ws2.Range("D13").Value = 4.6
$ws2.Range("E13").Value = 2
$ws2.Range("F13").Value = 23
$ws2.Range("G13").Value = 52
$ws2.Range("H13").Value = 5

$ws2.Range("A14").Value = "Sir Wheeze (Marcus Daley)"
$ws2.Range("B14").Value = 6.8
$ws2.Range("C14").Value = 0.4
$ws2.Range("D14").Value = 2.8
$ws2.Range("E14").Value = 2
$ws2.Range("F14").Value = 14
$ws2.Range("G14").Value = 34
$ws2.Range("H14").Value = 5

$ws2.Range("A15").Value = "Duke Hogs (Hogan Daley)"
$ws2.Range("B15").Value = 5.25
$ws2.Range("C15").Value = 0.25
$ws2.Range("D15").Value = 2.25
$ws2.Range("E15").Value = 1
$ws2.Range("F15").Value = 9
$ws2.Range("G15").Value = 21
$ws2.Range("H15").Value = 4

$ws2.Range("A16").Value = "Faith Youngquist"
$ws2.Range("B16").Value = 5
$ws2.Range("C16").Value = 0.3333333333333333
$ws2.Range("D16").Value = 2
$ws2.Range("E16").Value = 2
$ws2.Range("F16").Value = 12
$ws2.Range("G16").Value = 30
$ws2.Range("H16").Value = 6

# --- Sheet3: Team Rosters (content unchanged; refresh selection) ---
$ws3 = $wb.Worksheets.Item("Team Rosters")
$ws3.Range("A8").Select()

# Individual Stats ends up as the active/selected tab
$ws2.Activate()
$ws2.Range("I22").Select()

